$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("swiss_matches")
$ws2 = $wb.Worksheets.Item("teams")

# --- swiss_matches: drop the old "placement" rows (5-8), keep only the
# 4-row incidence/score matrix, and refresh the match data itself.
$ws1.Rows("5:8").Delete()

$ws1.Range("A1").Value = "Cosmo"
$ws1.Range("B1").Value = "Vertigo"
$ws1.Range("C1").Value = 15
$ws1.Range("D1").Value = 6

$ws1.Range("A2").Value = "Frizmi"
$ws1.Range("B2").Value = "Nuclear Discs"
$ws1.Range("C2").Value = 12
$ws1.Range("D2").Value = 15

$ws1.Range("A3").Value = "Frizmi"
$ws1.Range("B3").Value = "Cosmo"
$ws1.Range("C3").Value = 10
$ws1.Range("D3").Value = 6

$ws1.Range("A4").Value = "Vertigo"
$ws1.Range("B4").Value = "Nuclear Discs"
$ws1.Range("C4").Value = 15
$ws1.Range("D4").Value = 8

# --- teams: swap out the stale/unused team names for the current list.
$ws2.Range("A2").Value = "Frizmi"
$ws2.Range("A3").Value = "Vertigo"
$ws2.Range("A4").Value = "Nuclear Discs"
$ws2.Range("A5").ClearContents()

# --- view state: swiss_matches becomes the selected/active sheet.
$ws2.Range("A4").Select()
$ws1.Activate()
$ws1.Range("E9").Select()
